# Generate Report for Handoff
#
# The localization handoff run moved "low" priority items to "ht" (handed-off)
# and refreshed the handoff timestamps for the zh-cn and de-de sheets'
# "Ready for handoff" rows (rows 4-7).

$wb = $excel.ActiveWorkbook

$ws_zhcn = $wb.Sheets.Item("zh-cn")
$ws_dede = $wb.Sheets.Item("de-de")

# zh-cn: Priority low -> ht, Latest Handoff Datetime refreshed
foreach ($row in 4..7) {
    $ws_zhcn.Range("E" + $row).Value = "ht"
    $ws_zhcn.Range("H" + $row).Value = "2016-08-17 16:31:34"
}

# de-de: Priority low -> ht, Latest Handoff Datetime refreshed
foreach ($row in 4..7) {
    $ws_dede.Range("E" + $row).Value = "ht"
    $ws_dede.Range("H" + $row).Value = "2016-08-17 16:31:39"
}
